$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.214.87"
$ws.Range("E2").Value = "  +5.59%  "
$ws.Range("D3").Value = "1.785.30"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2681"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06276"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").Value = "1.780.15"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("E11").Value = "  +3.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07048"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6293"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.662"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.49%  "
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").Value = "28.183.96"
$ws.Range("E16").Value = "  +6.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9998"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9996"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007240"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +5.16%  "
$ws.Range("D21").Value = "2.008.55"
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.758"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.263"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.861"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "109.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.387"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.190"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08280"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.769"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04898"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.083"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.616"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6537"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9492"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.618"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.048"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.925"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01555"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9995"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.192"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1214"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05447"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.976"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("E50").Value = "  +4.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.08%  "
